$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.377.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.106.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.46%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5229'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.95%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4437'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.87%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.67'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09362'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.171'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.90'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.659'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.03%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.924'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.57%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.042.87'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '101.84'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001161'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.008'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.21'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06720'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.355'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.94%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.426.29'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.61'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.300'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.97'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.68'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.518'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '133.72'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.143'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.690'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.63%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.761'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +11.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.248'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.924'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.16%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02636'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06783'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.7034'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.351'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.85%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.55'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2224'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6849'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.51'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.354'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.47%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.399'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +20.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.643'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000350'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.212'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +9.21%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.11%  '
